# Edit script: apply the "Schedule and Slots" update to the workbook.
$wb = $excel.ActiveWorkbook

# --- Rename the "Slots" sheet to "Schedule and Slots" ---
$ws = $wb.Worksheets.Item("Slots")
$ws.Name = "Schedule and Slots"

# --- Shrink column D on the "Patients" sheet ---
$patients = $wb.Worksheets.Item("Patients")
$patients.Columns.Item(4).ColumnWidth = 34.59

# --- Build the new "Requirements for Schedule and Slots ..." block (rows 20-25) ---

# Row 20: merged title banner, bold, centered, grey fill
$r20 = $ws.Range("B20:F20")
$r20.Interior.ThemeColor = 2
$r20.Font.Bold = $true
$r20.Font.Size = 12
$r20.HorizontalAlignment = -4108
$r20.Merge()
$ws.Range("B20").Value = "Requirements for Schedule and Slots for Tests relating to serviceCategory and serviceType"

# Row 21: Schedule Day 1 / Schedule Should have NO serviceCategory set
$b21 = $ws.Range("B21")
$b21.Interior.ThemeColor = 2
$b21.HorizontalAlignment = -4131
$b21.Value = "Schedule Day 1"

$c21 = $ws.Range("C21:F21")
$c21.Interior.ThemeColor = 2
$c21.HorizontalAlignment = -4131
$c21.Merge()
$ws.Range("C21").Value = "Schedule Should have NO serviceCategory set"

# Row 23: Slots on Day 1 / Slots Should have NO serviceType set
$b23 = $ws.Range("B23")
$b23.Interior.ThemeColor = 2
$b23.HorizontalAlignment = -4131
$b23.Value = "Slots on Day 1"

$c23 = $ws.Range("C23:F23")
$c23.Interior.ThemeColor = 2
$c23.HorizontalAlignment = -4131
$c23.Merge()
$ws.Range("C23").Value = "Slots Should have NO serviceType set"

# Row 24: Slots on Day 2 / Slots Should have serviceType set to something
$b24 = $ws.Range("B24")
$b24.Interior.ThemeColor = 2
$b24.HorizontalAlignment = -4131
$b24.Value = "Slots on Day 2"

$c24 = $ws.Range("C24:F24")
$c24.Interior.ThemeColor = 2
$c24.HorizontalAlignment = -4131
$c24.Merge()
$ws.Range("C24").Value = "Slots Should have serviceType set to something"

# Row 22: Schedule Day 2 / Schedule Should have serviceCategory set to something
$b22 = $ws.Range("B22")
$b22.Interior.ThemeColor = 2
$b22.HorizontalAlignment = -4131
$b22.Value = "Schedule Day 2"

$c22 = $ws.Range("C22:F22")
$c22.Interior.ThemeColor = 2
$c22.HorizontalAlignment = -4131
$c22.Merge()
$ws.Range("C22").Value = "Schedule Should have serviceCategory set to something"

# Row 25: blank grey spacer row
$ws.Range("B25:F25").Interior.ThemeColor = 2

# --- Select / activate the "Schedule and Slots" sheet as the visible tab ---
$ws.Activate()
$ws.Range("B28").Select()

Write-Output "edit complete"
